# Adapt column header formatting to respective input file names.
# - Rename header cells from "<Name>_old" / "<Name>_new" suffixes to
#   "<Name>_FV2410" / "<Name>_FV2504" (the "diff" header stays as-is).
# - Turn the data range into a native Excel Table ("Table1").
# - Freeze the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410",
  "Segment ID_FV2410", "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410", "Bedingung_FV2410",
  "diff",
  "Segmentname_FV2504", "Segmentgruppe_FV2504", "Segment_FV2504", "Datenelement_FV2504",
  "Segment ID_FV2504", "Code_FV2504", "Qualifier_FV2504", "Beschreibung_FV2504",
  "Bedingungsausdruck_FV2504", "Bedingung_FV2504"
)

$headerRange = $ws.Range("A1:U1")

# Remember current header formatting (bold, fill, border, alignment) so it can be
# re-applied after the table is created without Excel baking a header-row dxf
# (which happens whenever ListObjects.Add() sees a non-default header style).
$headerRange.ClearFormats()

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$dataRange = $ws.Range("A1:U61")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Re-apply the original header-row formatting (bold font, grey fill, thin border,
# centered + wrapped text) now that the table already exists.
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 14277081
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true

# Freeze the header row.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
